# Spring-MVC-Project-Assignment.docx edit
#
# 1) In the "Code quality" bullet, the stray "_GoBack" bookmark that sat
#    between "following SOLID" and " principles" is removed and the two
#    runs are merged back into a single run "following SOLID principles".
# 2) Two new paragraphs are appended at the very end of the document body:
#    an empty one, followed by one containing "CHECK ALL TODOS!" - this
#    paragraph now carries the (re-inserted) "_GoBack" bookmark.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "following SOLID" + <bookmark _GoBack/> + " principles"
#         -> "following SOLID principles"  (bookmark dropped)
# ---------------------------------------------------------------------

# Locate the whole "Code quality (...)" bullet paragraph that holds the text.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Code quality", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$bulletPara = $d.Range($anchor.Start, $anchor.End)
$bulletPara.Expand(4) | Out-Null   # wdParagraph - grabs the full paragraph

# Drop the trailing paragraph mark from the delete range so the paragraph
# itself (and its pPr) stays put; only its run content is rebuilt.
$deleteRange = $d.Range($bulletPara.Start, $bulletPara.End - 1)
$deleteRange.Delete()

$rebuiltParagraphXml = @'
<w:p w:rsidR="0018280F" w:rsidRPr="005D080C" w:rsidRDefault="00F065A8" w:rsidP="00A0270C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="005D080C"><w:rPr><w:b/><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Code </w:t></w:r><w:r w:rsidR="00A0270C" w:rsidRPr="005D080C"><w:rPr><w:b/><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t>q</w:t></w:r><w:r w:rsidRPr="005D080C"><w:rPr><w:b/><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t>uality</w:t></w:r><w:r w:rsidR="003946AD" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0096761E" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">(well-structured code, </w:t></w:r><w:r w:rsidR="00EB5866" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">following </w:t></w:r><w:r w:rsidR="00A0270C" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">the </w:t></w:r><w:r w:rsidR="00EB5866" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t>MVC pattern</w:t></w:r><w:r w:rsidR="0096761E" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="005442BC" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t>following SOLID principles</w:t></w:r><w:r w:rsidR="0096761E" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, etc.) </w:t></w:r><w:r w:rsidR="003946AD" w:rsidRPr="005D080C"><w:rPr><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">– </w:t></w:r><w:r w:rsidRPr="005D080C"><w:rPr><w:b/><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t>0…</w:t></w:r><w:r w:rsidR="00904FBF" w:rsidRPr="005D080C"><w:rPr><w:b/><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r><w:r w:rsidR="00A25BEC" w:rsidRPr="005D080C"><w:rPr><w:b/><w:highlight w:val="green"/><w:lang w:val="en-US"/></w:rPr><w:t>0</w:t></w:r></w:p>
'@

$insertionPoint = $d.Range($deleteRange.Start, $deleteRange.Start)
$insertionPoint.InsertXML($rebuiltParagraphXml)

# ---------------------------------------------------------------------
# Edit 2: append an empty paragraph, then a "CHECK ALL TODOS!" paragraph
#         carrying the "_GoBack" bookmark, at the end of the document.
# ---------------------------------------------------------------------

$endOfDoc = $d.Content
$endOfDoc.Collapse(0) | Out-Null   # wdCollapseEnd

$newParagraphsXml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>CHECK ALL TODOS!</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>
'@

$endOfDoc.InsertXML($newParagraphsXml)

Write-Host "Applied: merged SOLID-principles run and appended CHECK ALL TODOS paragraph."
